$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.Style = "Normal"
}

Set-TextValue "D2" "308.80"
Set-TextValue "E2" "0.54%"
Set-TextValue "D3" "41.17"
Set-TextValue "E3" "4.52%"
Set-TextValue "D4" "5.125"
Set-TextValue "E4" "0.50%"
Set-TextValue "D5" "0.07644"
Set-TextValue "E5" "-0.66%"
Set-TextValue "D6" "1.621"
Set-TextValue "E6" "-0.94%"
Set-TextValue "E7" "0.89%"
Set-TextValue "D8" "0.9062"
Set-TextValue "E8" "-1.03%"
Set-TextValue "D9" "0.1125"
Set-TextValue "E9" "10.42%"
Set-TextValue "D10" "0.1802"
Set-TextValue "E10" "3.23%"
Set-TextValue "D11" "0.09107"
Set-TextValue "E11" "-2.54%"
Set-TextValue "D12" "0.04260"
Set-TextValue "E12" "-3.75%"
Set-TextValue "D13" "0.1051"
Set-TextValue "D14" "0.001259"
Set-TextValue "E14" "0.61%"
Set-TextValue "D15" "0.005828"
Set-TextValue "E15" "-0.73%"
Set-TextValue "E16" "-0.45%"
Set-TextValue "D17" "4.272"
Set-TextValue "E17" "0.62%"
Set-TextValue "E18" "0.37%"
Set-TextValue "D19" "6.737"
Set-TextValue "E19" "-3.54%"
Set-TextValue "D20" "0.1360"
Set-TextValue "E20" "0.92%"
Set-TextValue "E21" "-2.66%"
Set-TextValue "D22" "0.04064"
Set-TextValue "E22" "-1.70%"
Set-TextValue "D23" "0.001266"
Set-TextValue "E23" "5.53%"
Set-TextValue "E24" "-1.53%"
Set-TextValue "E25" "-2.13%"
Set-TextValue "D26" "0.0003747"
Set-TextValue "D38" "0.02423"
Set-TextValue "E38" "-1.47%"
Set-TextValue "D39" "0.05242"
Set-TextValue "E39" "1.00%"
Set-TextValue "D40" "0.007810"
Set-TextValue "E40" "-1.21%"
Set-TextValue "D41" "0.1301"
Set-TextValue "E41" "-1.43%"
Set-TextValue "D42" "0.006536"
Set-TextValue "E42" "-8.63%"
Set-TextValue "E43" "0.19%"
Set-TextValue "D44" "0.007581"
Set-TextValue "E44" "-9.47%"
Set-TextValue "D45" "0.3084"
Set-TextValue "E45" "0.78%"
Set-TextValue "D46" "0.00006777"
Set-TextValue "E46" "5.76%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.18%"
Set-TextValue "D48" "0.06286"
Set-TextValue "E48" "1,320.55%"
Set-TextValue "E49" "40.17%"
Set-TextValue "D50" "0.00002101"
Set-TextValue "E50" "0.18%"
Set-TextValue "D51" "0.0002001"
Set-TextValue "E51" "0.18%"
